$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# I2: expiry date, stored as literal text "2022-05-12" (not a real date value).
# Force text formatting first so Excel doesn't auto-convert the date-like
# string into a date serial number, then drop back to the Normal style so
# the cell doesn't keep a stray custom number format.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "2022-05-12"
$ws.Range("I2").Style = "Normal"

# J2: option instrument symbol
$ws.Range("J2").Value = "NIFTY2251216000CE"

# K2: stoploss_type Value -> Percentage
$ws.Range("K2").Value = "Percentage"

# O2: target_type Value -> Percentage
$ws.Range("O2").Value = "Percentage"
